$excel.UserName = "Author"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")
$ws.Range("B17").Value = "   Trim me!   "
$ws.Range("B18").Value = "                                   "
$comment = $ws.Range("B18").AddComment("OpenL User:" + [char]10 + "This cell contains spaces only.")
$ws.Range("C20").Select()
